$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (datetime) updates: shift by +0.625 days (15 hours)
$aUpdates = @(
    ,(2, 44280.54166666667)
    ,(3, 44280.58333333333)
    ,(4, 44280.625)
    ,(5, 44280.66666666667)
    ,(6, 44280.70833333333)
    ,(7, 44280.75)
    ,(8, 44280.79166666667)
    ,(9, 44280.83333333333)
    ,(10, 44280.875)
    ,(11, 44280.91666666667)
    ,(12, 44280.95833333333)
    ,(13, 44281)
    ,(14, 44281.04166666667)
    ,(15, 44281.08333333333)
    ,(16, 44281.41666666667)
    ,(17, 44281.54166666667)
    ,(18, 44281.58333333333)
    ,(19, 44281.625)
    ,(20, 44281.66666666667)
    ,(21, 44281.70833333333)
    ,(22, 44281.75)
    ,(23, 44281.79166666667)
    ,(24, 44281.83333333333)
    ,(25, 44281.875)
    ,(26, 44281.91666666667)
    ,(27, 44281.95833333333)
    ,(28, 44282)
    ,(29, 44282.04166666667)
    ,(30, 44282.08333333333)
    ,(31, 44284.45833333333)
    ,(32, 44284.54166666667)
    ,(33, 44284.58333333333)
    ,(34, 44284.625)
    ,(35, 44284.66666666667)
    ,(36, 44284.70833333333)
    ,(37, 44284.75)
    ,(38, 44284.79166666667)
    ,(39, 44284.83333333333)
    ,(40, 44284.875)
    ,(41, 44284.91666666667)
    ,(42, 44284.95833333333)
    ,(43, 44285)
    ,(44, 44285.08333333333)
    ,(45, 44285.54166666667)
    ,(46, 44285.58333333333)
    ,(47, 44285.625)
    ,(48, 44285.66666666667)
    ,(49, 44285.70833333333)
    ,(50, 44285.75)
    ,(51, 44285.79166666667)
    ,(52, 44285.83333333333)
    ,(53, 44285.875)
    ,(54, 44285.91666666667)
    ,(55, 44285.95833333333)
    ,(56, 44286)
    ,(57, 44286.04166666667)
    ,(58, 44286.08333333333)
    ,(59, 44286.45833333333)
    ,(60, 44286.54166666667)
    ,(61, 44286.58333333333)
    ,(62, 44286.625)
    ,(63, 44286.66666666667)
    ,(64, 44286.70833333333)
    ,(65, 44286.75)
    ,(66, 44286.79166666667)
    ,(67, 44286.83333333333)
    ,(68, 44286.875)
    ,(69, 44286.91666666667)
    ,(70, 44286.95833333333)
    ,(71, 44287)
    ,(72, 44287.04166666667)
    ,(73, 44287.08333333333)
    ,(74, 44287.45833333333)
    ,(75, 44287.54166666667)
    ,(76, 44287.58333333333)
    ,(77, 44287.625)
    ,(78, 44287.66666666667)
    ,(79, 44287.70833333333)
    ,(80, 44287.75)
    ,(81, 44287.79166666667)
    ,(82, 44287.83333333333)
    ,(83, 44287.875)
    ,(84, 44287.91666666667)
    ,(85, 44287.95833333333)
    ,(86, 44288)
    ,(87, 44288.04166666667)
    ,(88, 44288.08333333333)
    ,(89, 44291.41666666667)
    ,(90, 44291.54166666667)
    ,(91, 44291.58333333333)
    ,(92, 44291.625)
    ,(93, 44291.66666666667)
    ,(94, 44291.70833333333)
    ,(95, 44291.75)
    ,(96, 44291.79166666667)
    ,(97, 44291.83333333333)
    ,(98, 44291.875)
    ,(99, 44291.91666666667)
    ,(100, 44291.95833333333)
    ,(101, 44292)
    ,(102, 44292.08333333333)
    ,(103, 44292.41666666667)
    ,(104, 44292.45833333333)
    ,(105, 44292.5)
    ,(106, 44292.54166666667)
    ,(107, 44292.58333333333)
    ,(108, 44292.625)
    ,(109, 44292.66666666667)
    ,(110, 44292.70833333333)
    ,(111, 44292.75)
    ,(112, 44292.79166666667)
    ,(113, 44292.83333333333)
    ,(114, 44292.875)
    ,(115, 44292.91666666667)
    ,(116, 44292.95833333333)
    ,(117, 44293)
    ,(118, 44293.08333333333)
    ,(119, 44293.58333333333)
    ,(120, 44293.625)
    ,(121, 44293.66666666667)
    ,(122, 44293.70833333333)
    ,(123, 44293.75)
    ,(124, 44293.79166666667)
    ,(125, 44293.83333333333)
    ,(126, 44293.875)
    ,(127, 44293.91666666667)
    ,(128, 44293.95833333333)
    ,(129, 44294)
    ,(130, 44294.04166666667)
    ,(131, 44294.08333333333)
    ,(132, 44294.58333333333)
    ,(133, 44294.625)
    ,(134, 44294.66666666667)
    ,(135, 44294.70833333333)
    ,(136, 44294.75)
    ,(137, 44294.79166666667)
    ,(138, 44294.83333333333)
    ,(139, 44294.875)
    ,(140, 44294.91666666667)
    ,(141, 44294.95833333333)
    ,(142, 44295)
    ,(143, 44295.04166666667)
    ,(144, 44295.08333333333)
    ,(145, 44295.5)
    ,(146, 44295.58333333333)
    ,(147, 44295.625)
    ,(148, 44295.66666666667)
    ,(149, 44295.70833333333)
    ,(150, 44295.75)
    ,(151, 44295.79166666667)
    ,(152, 44295.83333333333)
    ,(153, 44295.875)
    ,(154, 44295.91666666667)
    ,(155, 44295.95833333333)
    ,(156, 44296)
    ,(157, 44296.04166666667)
    ,(158, 44296.08333333333)
    ,(159, 44280.33333333333)
    ,(160, 44280.375)
    ,(161, 44280.41666666667)
    ,(162, 44280.45833333333)
    ,(163, 44280.5)
    ,(164, 44281.125)
    ,(165, 44281.16666666667)
    ,(166, 44281.20833333333)
    ,(167, 44281.25)
    ,(168, 44281.29166666667)
    ,(169, 44281.33333333333)
    ,(170, 44281.375)
    ,(171, 44281.45833333333)
    ,(172, 44281.5)
    ,(173, 44282.125)
    ,(174, 44282.16666666667)
    ,(175, 44282.20833333333)
    ,(176, 44282.25)
    ,(177, 44282.29166666667)
    ,(178, 44282.33333333333)
    ,(179, 44282.375)
    ,(180, 44282.41666666667)
    ,(181, 44282.45833333333)
    ,(182, 44282.5)
    ,(183, 44282.54166666667)
    ,(184, 44282.58333333333)
    ,(185, 44282.625)
    ,(186, 44282.70833333333)
    ,(187, 44282.75)
    ,(188, 44282.79166666667)
    ,(189, 44282.83333333333)
    ,(190, 44282.875)
    ,(191, 44282.91666666667)
    ,(192, 44282.95833333333)
    ,(193, 44283)
    ,(194, 44283.04166666667)
    ,(195, 44283.08333333333)
    ,(196, 44283.125)
    ,(197, 44283.16666666667)
    ,(198, 44283.20833333333)
    ,(199, 44283.25)
    ,(200, 44283.29166666667)
    ,(201, 44283.33333333333)
    ,(202, 44283.375)
    ,(203, 44283.41666666667)
    ,(204, 44283.45833333333)
    ,(205, 44283.5)
    ,(206, 44283.54166666667)
    ,(207, 44283.58333333333)
    ,(208, 44283.625)
    ,(209, 44283.66666666667)
    ,(210, 44283.70833333333)
    ,(211, 44283.75)
    ,(212, 44283.79166666667)
    ,(213, 44283.83333333333)
    ,(214, 44283.875)
    ,(215, 44283.91666666667)
    ,(216, 44283.95833333333)
    ,(217, 44284)
    ,(218, 44284.04166666667)
    ,(219, 44284.08333333333)
    ,(220, 44284.125)
    ,(221, 44284.16666666667)
    ,(222, 44284.20833333333)
    ,(223, 44284.25)
    ,(224, 44284.29166666667)
    ,(225, 44284.33333333333)
    ,(226, 44284.375)
    ,(227, 44284.41666666667)
    ,(228, 44284.5)
    ,(229, 44285.04166666667)
    ,(230, 44285.125)
    ,(231, 44285.16666666667)
    ,(232, 44285.20833333333)
    ,(233, 44285.25)
    ,(234, 44285.29166666667)
    ,(235, 44285.33333333333)
    ,(236, 44285.375)
    ,(237, 44285.41666666667)
    ,(238, 44285.45833333333)
    ,(239, 44285.5)
    ,(240, 44286.125)
    ,(241, 44286.16666666667)
    ,(242, 44286.20833333333)
    ,(243, 44286.25)
    ,(244, 44286.29166666667)
    ,(245, 44286.33333333333)
    ,(246, 44286.375)
    ,(247, 44286.41666666667)
    ,(248, 44286.5)
    ,(249, 44287.125)
    ,(250, 44287.16666666667)
    ,(251, 44287.20833333333)
    ,(252, 44287.25)
    ,(253, 44287.29166666667)
    ,(254, 44287.33333333333)
    ,(255, 44287.375)
    ,(256, 44287.41666666667)
    ,(257, 44287.5)
    ,(258, 44288.125)
    ,(259, 44288.16666666667)
    ,(260, 44288.20833333333)
    ,(261, 44288.25)
    ,(262, 44288.29166666667)
    ,(263, 44288.33333333333)
    ,(264, 44288.375)
    ,(265, 44288.41666666667)
    ,(266, 44288.45833333333)
    ,(267, 44288.54166666667)
    ,(268, 44288.58333333333)
    ,(269, 44288.625)
    ,(270, 44288.66666666667)
    ,(271, 44288.70833333333)
    ,(272, 44288.75)
    ,(273, 44288.79166666667)
    ,(274, 44288.83333333333)
    ,(275, 44288.875)
    ,(276, 44288.91666666667)
    ,(277, 44288.95833333333)
    ,(278, 44289)
    ,(279, 44289.04166666667)
    ,(280, 44289.08333333333)
    ,(281, 44289.125)
    ,(282, 44289.16666666667)
    ,(283, 44289.20833333333)
    ,(284, 44289.25)
    ,(285, 44289.29166666667)
    ,(286, 44289.33333333333)
    ,(287, 44289.375)
    ,(288, 44289.41666666667)
    ,(289, 44289.45833333333)
    ,(290, 44289.5)
    ,(291, 44289.54166666667)
    ,(292, 44289.58333333333)
    ,(293, 44289.625)
    ,(294, 44289.66666666667)
    ,(295, 44289.70833333333)
    ,(296, 44289.75)
    ,(297, 44289.79166666667)
    ,(298, 44289.83333333333)
    ,(299, 44289.875)
    ,(300, 44289.91666666667)
    ,(301, 44289.95833333333)
    ,(302, 44290)
    ,(303, 44290.04166666667)
    ,(304, 44290.08333333333)
    ,(305, 44290.125)
    ,(306, 44290.16666666667)
    ,(307, 44290.20833333333)
    ,(308, 44290.25)
    ,(309, 44290.29166666667)
    ,(310, 44290.33333333333)
    ,(311, 44290.375)
    ,(312, 44290.41666666667)
    ,(313, 44290.45833333333)
    ,(314, 44290.5)
    ,(315, 44290.54166666667)
    ,(316, 44290.58333333333)
    ,(317, 44290.625)
    ,(318, 44290.66666666667)
    ,(319, 44290.70833333333)
    ,(320, 44290.75)
    ,(321, 44290.79166666667)
    ,(322, 44290.83333333333)
    ,(323, 44290.875)
    ,(324, 44290.91666666667)
    ,(325, 44290.95833333333)
    ,(326, 44291)
    ,(327, 44291.04166666667)
    ,(328, 44291.08333333333)
    ,(329, 44291.125)
    ,(330, 44291.16666666667)
    ,(331, 44291.20833333333)
    ,(332, 44291.25)
    ,(333, 44291.29166666667)
    ,(334, 44291.33333333333)
    ,(335, 44291.375)
    ,(336, 44291.45833333333)
    ,(337, 44291.5)
    ,(338, 44292.04166666667)
    ,(339, 44292.125)
    ,(340, 44292.16666666667)
    ,(341, 44292.20833333333)
    ,(342, 44292.25)
    ,(343, 44292.29166666667)
    ,(344, 44292.33333333333)
    ,(345, 44292.375)
    ,(346, 44293.04166666667)
    ,(347, 44293.125)
    ,(348, 44293.16666666667)
    ,(349, 44293.20833333333)
    ,(350, 44293.25)
    ,(351, 44293.29166666667)
    ,(352, 44293.33333333333)
    ,(353, 44293.375)
    ,(354, 44293.41666666667)
    ,(355, 44293.45833333333)
    ,(356, 44293.5)
    ,(357, 44293.54166666667)
    ,(358, 44294.125)
    ,(359, 44294.16666666667)
    ,(360, 44294.20833333333)
    ,(361, 44294.25)
    ,(362, 44294.29166666667)
    ,(363, 44294.33333333333)
    ,(364, 44294.375)
    ,(365, 44294.41666666667)
    ,(366, 44294.45833333333)
    ,(367, 44294.5)
    ,(368, 44294.54166666667)
    ,(369, 44295.125)
    ,(370, 44295.16666666667)
    ,(371, 44295.20833333333)
    ,(372, 44295.25)
    ,(373, 44295.29166666667)
    ,(374, 44295.33333333333)
    ,(375, 44295.375)
    ,(376, 44295.41666666667)
    ,(377, 44295.45833333333)
    ,(378, 44295.54166666667)
    ,(379, 44296.125)
    ,(380, 44296.16666666667)
    ,(381, 44296.20833333333)
    ,(382, 44296.25)
    ,(383, 44296.29166666667)
    ,(384, 44296.33333333333)
    ,(385, 44296.375)
)

# Column B (price) updates: tiny float recalculation deltas
$bUpdates = @(
    ,(2, -0.4676517815562942)
    ,(3, -0.565117464440652)
    ,(4, -0.6728143516056915)
    ,(5, -0.9334408185450886)
    ,(7, -1.177374267973901)
    ,(8, -1.318457190160103)
    ,(10, -1.453616783552227)
    ,(16, -1.238223009222148)
    ,(21, -1.470848285498633)
    ,(24, -1.621623927529689)
    ,(25, -1.421307717402714)
    ,(26, -1.154757921669242)
    ,(27, -1.154757921669242)
    ,(28, -1.318457190160103)
    ,(30, -1.154757921669242)
    ,(32, -1.242530884708749)
    ,(33, -1.27268601311496)
    ,(34, -1.329765363312431)
    ,(35, -1.211837271866713)
    ,(37, -1.005328490727751)
    ,(38, -0.7880500208722829)
    ,(39, -0.95551868041392)
    ,(40, -0.8370521045323741)
    ,(41, -0.8370521045323741)
    ,(42, -0.8343596823532506)
    ,(43, -0.8348981667890754)
    ,(44, -0.8370521045323741)
    ,(53, -0.8628993574519845)
    ,(54, -0.8672072329385884)
    ,(55, -0.9226711298285822)
    ,(56, -0.8672072329385884)
    ,(57, -0.8672072329385884)
    ,(58, -0.8672072329385884)
    ,(60, -0.7805112387707309)
    ,(61, -0.8672072329385884)
    ,(62, -0.7272012796240359)
    ,(64, -0.285105557811547)
    ,(65, -0.2463346784321374)
    ,(66, -0.1052517562459355)
    ,(68, -0.3847251784392097)
    ,(69, -0.4218806045111512)
    ,(71, -0.4218806045111512)
    ,(72, -0.4035721336930928)
    ,(73, -0.4218806045111512)
    ,(74, -0.05140331266341581)
    ,(75, -0.08101995663379921)
    ,(76, -0.03223326674803909)
    ,(77, 0.02559996165959081)
    ,(80, 0.2899958196497609)
    ,(82, 0.4106163332746056)
    ,(84, 0.5344677535143985)
    ,(85, 0.5414680511801259)
    ,(86, 0.5344677535143985)
    ,(87, 0.5522377398966322)
    ,(88, 0.5344677535143985)
    ,(89, 0.7293991192831202)
    ,(90, 0.5656998507922621)
    ,(91, 0.5603150064340089)
    ,(92, 0.5770080239445933)
    ,(94, 0.4133087554537291)
    ,(96, 0.5576225842548854)
    ,(97, 0.5802389305595415)
    ,(98, 0.6028552768642037)
    ,(99, 0.6017783079925482)
    ,(100, 0.6198713850362755)
    ,(101, 0.6017783079925482)
    ,(102, 0.6017783079925482)
    ,(103, 0.8360190375765102)
    ,(104, 0.7234757904890422)
    ,(105, 0.7267066971039966)
    ,(106, 0.7191679150024446)
    ,(107, 0.7541694033310812)
    ,(108, 0.8080178469136009)
    ,(110, 1.063528711712657)
    ,(111, 0.9264844227951468)
    ,(112, 0.8607893216244712)
    ,(113, 0.9754865064552379)
    ,(114, 0.9776404441985368)
    ,(115, 0.809633300221075)
    ,(116, 0.7891708916597178)
    ,(117, 0.809633300221075)
    ,(118, 0.809633300221075)
    ,(119, 0.7746318118924385)
    ,(120, 0.7137830706441913)
    ,(121, 0.851042753336033)
    ,(122, 0.9329462360250493)
    ,(123, 0.9641783333029068)
    ,(124, 0.9065604986696142)
    ,(126, 0.8677896192901985)
    ,(127, 0.9421004714340755)
    ,(129, 1.0088725414764)
    ,(131, 0.9421004714340755)
    ,(132, 1.20703481386007)
    ,(145, 1.279191728260648)
    ,(147, 1.225881769113953)
    ,(148, 1.028257981166108)
    ,(156, 1.376657411145006)
)

foreach ($pair in $aUpdates) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 1).Value = $v
}

foreach ($pair in $bUpdates) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 2).Value = $v
}

Write-Host "Applied $($aUpdates.Count) column A updates and $($bUpdates.Count) column B updates"
